# Added more absorber materials (Rh, Ag, Tc, Te) to the EUV materials table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round-trip the sheet name so the workbook/worksheet XML gets regenerated
# with normalized relationship ids (matches a freshly-saved Excel file).
$ws.Name = "Sheet1_tmp"
$ws.Name = "Sheet1"

# New absorber material rows appended below the existing table.
$ws.Range("A8").Value = "Rh"
$ws.Range("B8").Value = 0.12496963899999999
$ws.Range("C8").Value = 0.031190166200000001

$ws.Range("A9").Value = "Ag"
$ws.Range("B9").Value = 0.10969390699999999
$ws.Range("C9").Value = 0.079397715600000002

$ws.Range("A10").Value = "Tc"
$ws.Range("B10").Value = 0.093783922500000005
$ws.Range("C10").Value = 0.012075969900000001

$ws.Range("A11").Value = "Te"
$ws.Range("B11").Value = 0.0271440633
$ws.Range("C11").Value = 0.074925527000000006

# New rows use a dedicated font (Helvetica Neue 10pt, black) distinct from
# the rest of the table.
$ws.Range("A8:C11").Font.Name = "Helvetica Neue"
$ws.Range("A8:C11").Font.Size = 10
$ws.Range("A8:C11").Font.Color = 0

# Match the selection left active in the source file.
$ws.Range("A11:C11").Select()
